$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.067937038739814
$ws.Range("D2").Value = 7.775327537238653
$ws.Range("E2").Value = 12.79089778093609
$ws.Range("F2").Value = 43.55292226603628
$ws.Range("G2").Value = 55.81828364157815
$ws.Range("H2").Value = 19.41613829829462
$ws.Range("I2").Value = 30.49359561336957
$ws.Range("J2").Value = 10.09421402311757
$ws.Range("C3").Value = 4.079197275270284
$ws.Range("D3").Value = 7.783950329519477
$ws.Range("E3").Value = 12.7962115933908
$ws.Range("F3").Value = 42.80830285766458
$ws.Range("G3").Value = 54.23993323042035
$ws.Range("H3").Value = 19.23223933023537
$ws.Range("I3").Value = 29.98868480780469
$ws.Range("J3").Value = 10.09965827486539
$ws.Range("C4").Value = 4.086823181783126
$ws.Range("D4").Value = 7.789593014885036
$ws.Range("E4").Value = 12.80222834265474
$ws.Range("F4").Value = 42.36012949402359
$ws.Range("G4").Value = 53.26890903784906
$ws.Range("H4").Value = 19.12460663806862
$ws.Range("I4").Value = 29.68522153962714
$ws.Range("J4").Value = 10.10513802733398
$ws.Range("C5").Value = 4.090109392622248
$ws.Range("D5").Value = 7.791980221879072
$ws.Range("E5").Value = 12.80537095394674
$ws.Range("F5").Value = 42.18000477163207
$ws.Range("G5").Value = 52.87339574061149
$ws.Range("H5").Value = 19.08211190537824
$ws.Range("I5").Value = 29.56336527187999
$ws.Range("J5").Value = 10.10790694008688
$ws.Range("C6").Value = 4.090665838468226
$ws.Range("D6").Value = 7.792381922476789
$ws.Range("E6").Value = 12.80593444346414
$ws.Range("F6").Value = 42.15025368690453
$ws.Range("G6").Value = 52.80775238547449
$ws.Range("H6").Value = 19.07513927552349
$ws.Range("I6").Value = 29.54324488440706
$ws.Range("J6").Value = 10.10839903397143
$ws.Range("C7").Value = 4.086866778263603
$ws.Range("D7").Value = 7.789624853964197
$ws.Range("E7").Value = 12.80226793078253
$ws.Range("F7").Value = 42.35768979431256
$ws.Range("G7").Value = 53.26357334810401
$ws.Range("H7").Value = 19.12402795819506
$ws.Range("I7").Value = 29.68357061698928
$ws.Range("J7").Value = 10.10517320223398
$ws.Range("C8").Value = 4.071671495189776
$ws.Range("D8").Value = 7.778228512134253
$ws.Range("E8").Value = 12.79215737886159
$ws.Range("F8").Value = 43.29444933871008
$ws.Range("G8").Value = 55.27485639963601
$ws.Range("H8").Value = 19.35165752383276
$ws.Range("I8").Value = 30.31823872306765
$ws.Range("J8").Value = 10.09564683369204
$ws.Range("C9").Value = 4.047544973677834
$ws.Range("D9").Value = 7.758634923494748
$ws.Range("E9").Value = 12.79425332100533
$ws.Range("F9").Value = 45.1923579334902
$ws.Range("G9").Value = 59.17608532848482
$ws.Range("H9").Value = 19.83812848513057
$ws.Range("I9").Value = 31.60769287801838
$ws.Range("J9").Value = 10.09398137337356
$ws.Range("C10").Value = 4.033306825936791
$ws.Range("D10").Value = 7.745907050268944
$ws.Range("E10").Value = 12.80924056609112
$ws.Range("F10").Value = 46.60957261525871
$ws.Range("G10").Value = 61.981652212655
$ws.Range("H10").Value = 20.21737083869697
$ws.Range("I10").Value = 32.57281943039924
$ws.Range("J10").Value = 10.1032033125207
$ws.Range("C11").Value = 4.027593510409064
$ws.Range("D11").Value = 7.740476480735375
$ws.Range("E11").Value = 12.81899056674694
$ws.Range("F11").Value = 47.25628024516151
$ws.Range("G11").Value = 63.23832606388071
$ws.Range("H11").Value = 20.39404390176299
$ws.Range("I11").Value = 33.01373424058502
$ws.Range("J11").Value = 10.1096781573947
$ws.Range("C12").Value = 4.025540408110646
$ws.Range("D12").Value = 7.738471567625717
$ws.Range("E12").Value = 12.82310469558071
$ws.Range("F12").Value = 47.50122072253101
$ws.Range("G12").Value = 63.71089757901819
$ws.Range("H12").Value = 20.4614908354891
$ws.Range("I12").Value = 33.18080441632898
$ws.Range("J12").Value = 10.11245835723054
$ws.Range("C13").Value = 4.025977660591156
$ws.Range("D13").Value = 7.738901072355536
$ws.Range("E13").Value = 12.82219987252577
$ws.Range("F13").Value = 47.4484701162485
$ws.Range("G13").Value = 63.60927505002352
$ws.Range("H13").Value = 20.44694155215289
$ws.Range("I13").Value = 33.14482073616267
$ws.Range("J13").Value = 10.11184498395784
$ws.Range("C14").Value = 4.027422384339081
$ws.Range("D14").Value = 7.740310503551593
$ws.Range("E14").Value = 12.81932057978253
$ws.Range("F14").Value = 47.27643207546961
$ws.Range("G14").Value = 63.27727364430339
$ws.Range("H14").Value = 20.39958218944831
$ws.Range("I14").Value = 33.02747801478515
$ws.Range("J14").Value = 10.10990030410755
$ws.Range("C15").Value = 4.028321715129676
$ws.Range("D15").Value = 7.741180526475396
$ws.Range("E15").Value = 12.81761189429881
$ws.Range("F15").Value = 47.17105279049238
$ws.Range("G15").Value = 63.07346929284409
$ws.Range("H15").Value = 20.37064257097947
$ws.Range("I15").Value = 32.95561114459593
$ws.Range("J15").Value = 10.10875189896313
$ws.Range("C16").Value = 4.033695622038209
$ws.Range("D16").Value = 7.746269169429702
$ws.Range("E16").Value = 12.80866242263742
$ws.Range("F16").Value = 46.56732825125148
$ws.Range("G16").Value = 61.89908879863986
$ws.Range("H16").Value = 20.20590394254974
$ws.Range("I16").Value = 32.54402816467346
$ws.Range("J16").Value = 10.10282608747748
$ws.Range("C17").Value = 4.037188354973043
$ws.Range("D17").Value = 7.749482826417786
$ws.Range("E17").Value = 12.8039235813655
$ws.Range("F17").Value = 46.19730219973535
$ws.Range("G17").Value = 61.17326632241945
$ws.Range("H17").Value = 20.10586883792651
$ws.Range("I17").Value = 32.29189684943348
$ws.Range("J17").Value = 10.09977510984827
$ws.Range("C18").Value = 4.039269128900386
$ws.Range("D18").Value = 7.751365074753489
$ws.Range("E18").Value = 12.80147391798407
$ws.Range("F18").Value = 45.98467618138523
$ws.Range("G18").Value = 60.75397454483058
$ws.Range("H18").Value = 20.04872548346197
$ws.Range("I18").Value = 32.14706349118176
$ws.Range("J18").Value = 10.09823483244165
$ws.Range("C19").Value = 4.039985965060445
$ws.Range("D19").Value = 7.752008188215201
$ws.Range("E19").Value = 12.80069188997181
$ws.Range("F19").Value = 45.91272709304165
$ws.Range("G19").Value = 60.61171286734235
$ws.Range("H19").Value = 20.02944701854534
$ws.Range("I19").Value = 32.09806242153542
$ws.Range("J19").Value = 10.09775015397212
$ws.Range("C20").Value = 4.036809107536698
$ws.Range("D20").Value = 7.749137226204418
$ws.Range("E20").Value = 12.8043994713039
$ws.Range("F20").Value = 46.23667286139261
$ws.Range("G20").Value = 61.25072317052071
$ws.Range("H20").Value = 20.11647732707277
$ws.Range("I20").Value = 32.31871862397818
$ws.Range("J20").Value = 10.1000776807827
$ws.Range("C21").Value = 4.026995032819068
$ws.Range("D21").Value = 7.73989512227204
$ws.Range("E21").Value = 12.82015484309045
$ws.Range("F21").Value = 47.32696445070253
$ws.Range("G21").Value = 63.37488382026311
$ws.Range("H21").Value = 20.41347842282201
$ws.Range("I21").Value = 33.06194283953453
$ws.Range("J21").Value = 10.11046259045549
$ws.Range("C22").Value = 4.021224779697756
$ws.Range("D22").Value = 7.734155135374947
$ws.Range("E22").Value = 12.83291167125691
$ws.Range("F22").Value = 48.03969709016926
$ws.Range("G22").Value = 64.74369508337084
$ws.Range("H22").Value = 20.61073578615554
$ws.Range("I22").Value = 33.54822457537134
$ws.Range("J22").Value = 10.11916348416579
$ws.Range("C23").Value = 4.024245366627255
$ws.Range("D23").Value = 7.737191250719672
$ws.Range("E23").Value = 12.82587799992593
$ws.Range("F23").Value = 47.65936106480515
$ws.Range("G23").Value = 64.01506375992504
$ws.Range("H23").Value = 20.5051851481821
$ws.Range("I23").Value = 33.28869005859334
$ws.Range("J23").Value = 10.11434443829059
$ws.Range("C24").Value = 4.03698033873114
$ws.Range("D24").Value = 7.749293364047076
$ws.Range("E24").Value = 12.80418346567879
$ws.Range("F24").Value = 46.21887303381959
$ws.Range("G24").Value = 61.21571114275333
$ws.Range("H24").Value = 20.11168007841102
$ws.Range("I24").Value = 32.30659211250418
$ws.Range("J24").Value = 10.09994022265615
$ws.Range("C25").Value = 4.053461589684801
$ws.Range("D25").Value = 7.763641839373895
$ws.Range("E25").Value = 12.79132965375957
$ws.Range("F25").Value = 44.67388760418544
$ws.Range("G25").Value = 58.12901065079848
$ws.Range("H25").Value = 19.70249413479836
$ws.Range("I25").Value = 31.25505110455875
$ws.Range("J25").Value = 10.09260199049064
